# Update the "Förändrad" (Changed) date column (C) for rows 2-12
# from serial date 45175 (2023-09-06) to 45177 (2023-09-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45177
    }
}
